# chore: adapt column header formatting to respective input file names
#
# The AHB-diff export used to label the "before" / "after" columns with a
# generic "_old" / "_new" suffix. This now uses the concrete format-version
# of each input file instead ("_FV2404" for the older, "_FV2410" for the
# newer formatversion) as the column-header suffix, turns the header +
# data range into a real Excel Table, and freezes the header row so it
# stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells -------------------------------------------
# "<name>_old" -> "<name>_FV2404" (older formatversion)
# "<name>_new" -> "<name>_FV2410" (newer formatversion)
[void]$ws.Cells.Replace("_old", "_FV2404")
[void]$ws.Cells.Replace("_new", "_FV2410")

# --- 2. Turn the header + data range into a real Excel Table -------------
$rng = $ws.Range("A1:U80")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row so it stays visible while scrolling --------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
